$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "97×47=4559"; New = "67×33=2211" },
    @{ Old = "37×45=1665"; New = "44×43=1892" },
    @{ Old = "43×49=2107"; New = "21×84=1764" },
    @{ Old = "34×48=1632"; New = "33×71=2343" },
    @{ Old = "18×57=1026"; New = "48×94=4512" },
    @{ Old = "75×49=3675"; New = "52×94=4888" },
    @{ Old = "20×19=380";  New = "75×50=3750" },
    @{ Old = "51×70=3570"; New = "15×64=960" },
    @{ Old = "86×34=2924"; New = "72×68=4896" },
    @{ Old = "54×52=2808"; New = "45×54=2430" },
    @{ Old = "56×24=1344"; New = "54×24=1296" },
    @{ Old = "59×87=5133"; New = "55×41=2255" },
    @{ Old = "73×56=4088"; New = "38×65=2470" },
    @{ Old = "61×37=2257"; New = "55×34=1870" },
    @{ Old = "52×41=2132"; New = "88×88=7744" },
    @{ Old = "69×26=1794"; New = "52×27=1404" },
    @{ Old = "47×70=3290"; New = "65×74=4810" },
    @{ Old = "33×34=1122"; New = "85×67=5695" },
    @{ Old = "78×19=1482"; New = "48×25=1200" },
    @{ Old = "69×69=4761"; New = "47×81=3807" },
    @{ Old = "23×92=2116"; New = "35×25=875" },
    @{ Old = "86×56=4816"; New = "90×81=7290" },
    @{ Old = "14×96=1344"; New = "30×17=510" },
    @{ Old = "27×51=1377"; New = "67×38=2546" },
    @{ Old = "40×79=3160"; New = "86×93=7998" }
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.Old, $true, $true, $false, $false, $false, $true, 1, $false, $r.New, 2)
}

$d.Save()
